$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up the data-output tags: hyphens -> underscores.
$ws.Range("A30").Value = "B-SOC_EMP_TNM"
$ws.Range("A31").Value = "B-SOC_EMP_TTN"

# Update the saved view: selection moved from B35 to A34, scrolled so row 11 / col A is top-left.
$ws.Activate()
$ws.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
